# "Generate Report for Handback" -- refresh the localization-status report
# after a successful handback: update each locale's Status / Latest
# Handback DateTime, clear the stale "handback not latest" Error Detail,
# and widen/narrow the affected columns to fit their new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- zh-cn ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("K2").Value = "2016-08-21 06:55:12"
$ws.Range("P2").Value = ""

$ws.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("K2").Value = "2016-08-21 06:55:18"
$ws.Range("P2").Value = ""

$ws.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- Overview --------------------------------------------------------------
# E2/F2 mirror the zh-cn/de-de Status cells (same shared text).
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus

$ws.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws.Columns.Item(6).ColumnWidth = 29.166666666666668
